# Fruta / hortaliza, semanal
#
# Insert one new weekly price record for "Feria Lagunitas de Puerto Montt - Mango"
# at row 260, pushing the existing rows 260:289 down to 261:290.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 260:289 down to 261:290, leaving a blank row 260 that inherits
# the formatting (incl. the date-formatted column D) of the row above it.
$ws.Rows.Item(260).Insert()

$ws.Range("A260").Value = 4
$ws.Range("B260").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C260").Value = "Los Lagos"
$ws.Range("D260").Value = 44918
$ws.Range("E260").Value = 10
$ws.Range("F260").Value = "Fruta"
$ws.Range("G260").Value = 100108
$ws.Range("H260").Value = "Tropicales y subtropicales"
$ws.Range("I260").Value = 100108002
$ws.Range("J260").Value = "Mango"
$ws.Range("K260").Value = "Sin especificar"
$ws.Range("L260").Value = "Primera"
$ws.Range("M260").Value = 300
$ws.Range("N260").Value = 8000
$ws.Range("O260").Value = 8500
$ws.Range("P260").Value = 8250
$ws.Range("Q260").Value = "$/bandeja 4 kilos"
$ws.Range("R260").Value = "Brasil"
$ws.Range("S260").Value = 2062
$ws.Range("T260").Value = 4
